$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.758258666666666
$ws.Cells.Item(2, 8).Value = 5.274775999999999
$ws.Cells.Item(2, 9).Value = 0.1132279568112417
$ws.Cells.Item(2, 10).Value = 0.1132279568112417
$ws.Cells.Item(2, 13).Value = 73.202511
$ws.Cells.Item(2, 14).Value = 219.607533
$ws.Cells.Item(2, 15).Value = 0.3264904632507938
$ws.Cells.Item(2, 16).Value = 0.3264904632507938
$ws.Cells.Item(2, 17).Value = 128.708949387512
$ws.Cells.Item(2, 18).Value = 1158.380544487608
$ws.Cells.Item(2, 19).Value = 0.03696784807224317
$ws.Cells.Item(2, 20).Value = 0.03696784807224317
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.758258666666666
$ws.Cells.Item(3, 8).Value = 5.274775999999999
$ws.Cells.Item(3, 9).Value = 0.1132279568112417
$ws.Cells.Item(3, 10).Value = 0.1132279568112417
$ws.Cells.Item(3, 15).Value = 0.4449719839907295
$ws.Cells.Item(3, 16).Value = 0.4449719839907295
$ws.Cells.Item(3, 17).Value = 175.4166905706222
$ws.Cells.Item(3, 18).Value = 1578.7502151356
$ws.Cells.Item(3, 19).Value = 0.05038326858551485
$ws.Cells.Item(3, 20).Value = 0.05038326858551485
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.758258666666666
$ws.Cells.Item(4, 8).Value = 5.274775999999999
$ws.Cells.Item(4, 9).Value = 0.1132279568112417
$ws.Cells.Item(4, 10).Value = 0.1132279568112417
$ws.Cells.Item(4, 13).Value = 39.54025133333334
$ws.Cells.Item(4, 14).Value = 118.620754
$ws.Cells.Item(4, 15).Value = 0.1763534446908907
$ws.Cells.Item(4, 16).Value = 0.1763534446908907
$ws.Cells.Item(4, 17).Value = 69.52198958901155
$ws.Cells.Item(4, 18).Value = 625.6979063011039
$ws.Cells.Item(4, 19).Value = 0.01996814021897387
$ws.Cells.Item(4, 20).Value = 0.01996814021897387
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.758258666666666
$ws.Cells.Item(5, 8).Value = 5.274775999999999
$ws.Cells.Item(5, 9).Value = 0.1132279568112417
$ws.Cells.Item(5, 10).Value = 0.1132279568112417
$ws.Cells.Item(5, 13).Value = 11.70021233333333
$ws.Cells.Item(5, 14).Value = 35.100637
$ws.Cells.Item(5, 15).Value = 0.05218410806758597
$ws.Cells.Item(5, 16).Value = 0.05218410806758598
$ws.Cells.Item(5, 17).Value = 20.57199973692355
$ws.Cells.Item(5, 18).Value = 185.147997632312
$ws.Cells.Item(5, 19).Value = 0.005908699934509792
$ws.Cells.Item(5, 20).Value = 0.005908699934509793
$ws.Cells.Item(6, 9).Value = 0.1732550390834427
$ws.Cells.Item(6, 10).Value = 0.1732550390834427
$ws.Cells.Item(6, 13).Value = 73.202511
$ws.Cells.Item(6, 14).Value = 219.607533
$ws.Cells.Item(6, 15).Value = 0.3264904632507938
$ws.Cells.Item(6, 16).Value = 0.3264904632507938
$ws.Cells.Item(6, 17).Value = 196.943181565105
$ws.Cells.Item(6, 18).Value = 1772.488634085945
$ws.Cells.Item(6, 19).Value = 0.05656611797088759
$ws.Cells.Item(6, 20).Value = 0.0565661179708876
$ws.Cells.Item(7, 9).Value = 0.1732550390834427
$ws.Cells.Item(7, 10).Value = 0.1732550390834427
$ws.Cells.Item(7, 15).Value = 0.4449719839907295
$ws.Cells.Item(7, 16).Value = 0.4449719839907295
$ws.Cells.Item(7, 19).Value = 0.07709363847735087
$ws.Cells.Item(7, 20).Value = 0.07709363847735087
$ws.Cells.Item(8, 9).Value = 0.1732550390834427
$ws.Cells.Item(8, 10).Value = 0.1732550390834427
$ws.Cells.Item(8, 13).Value = 39.54025133333334
$ws.Cells.Item(8, 14).Value = 118.620754
$ws.Cells.Item(8, 15).Value = 0.1763534446908907
$ws.Cells.Item(8, 16).Value = 0.1763534446908907
$ws.Cells.Item(8, 17).Value = 106.3786308842678
$ws.Cells.Item(8, 18).Value = 957.4076779584101
$ws.Cells.Item(8, 19).Value = 0.03055412295242002
$ws.Cells.Item(8, 20).Value = 0.03055412295242002
$ws.Cells.Item(9, 9).Value = 0.1732550390834427
$ws.Cells.Item(9, 10).Value = 0.1732550390834427
$ws.Cells.Item(9, 13).Value = 11.70021233333333
$ws.Cells.Item(9, 14).Value = 35.100637
$ws.Cells.Item(9, 15).Value = 0.05218410806758597
$ws.Cells.Item(9, 16).Value = 0.05218410806758598
$ws.Cells.Item(9, 17).Value = 31.47811475912278
$ws.Cells.Item(9, 18).Value = 283.303032832105
$ws.Cells.Item(9, 19).Value = 0.009041159682784204
$ws.Cells.Item(9, 20).Value = 0.009041159682784206
$ws.Cells.Item(10, 7).Value = 9.938311666666667
$ws.Cells.Item(10, 8).Value = 29.814935
$ws.Cells.Item(10, 9).Value = 0.6400052196548212
$ws.Cells.Item(10, 10).Value = 0.640005219654821
$ws.Cells.Item(10, 13).Value = 73.202511
$ws.Cells.Item(10, 14).Value = 219.607533
$ws.Cells.Item(10, 15).Value = 0.3264904632507938
$ws.Cells.Item(10, 16).Value = 0.3264904632507938
$ws.Cells.Item(10, 17).Value = 727.5093691005951
$ws.Cells.Item(10, 18).Value = 6547.584321905355
$ws.Cells.Item(10, 19).Value = 0.2089556006480286
$ws.Cells.Item(10, 20).Value = 0.2089556006480286
$ws.Cells.Item(11, 7).Value = 9.938311666666667
$ws.Cells.Item(11, 8).Value = 29.814935
$ws.Cells.Item(11, 9).Value = 0.6400052196548212
$ws.Cells.Item(11, 10).Value = 0.640005219654821
$ws.Cells.Item(11, 15).Value = 0.4449719839907295
$ws.Cells.Item(11, 16).Value = 0.4449719839907295
$ws.Cells.Item(11, 17).Value = 991.5183559033056
$ws.Cells.Item(11, 18).Value = 8923.665203129751
$ws.Cells.Item(11, 19).Value = 0.2847843923542284
$ws.Cells.Item(11, 20).Value = 0.2847843923542284
$ws.Cells.Item(12, 7).Value = 9.938311666666667
$ws.Cells.Item(12, 8).Value = 29.814935
$ws.Cells.Item(12, 9).Value = 0.6400052196548212
$ws.Cells.Item(12, 10).Value = 0.640005219654821
$ws.Cells.Item(12, 13).Value = 39.54025133333334
$ws.Cells.Item(12, 14).Value = 118.620754
$ws.Cells.Item(12, 15).Value = 0.1763534446908907
$ws.Cells.Item(12, 16).Value = 0.1763534446908907
$ws.Cells.Item(12, 17).Value = 392.9633411289989
$ws.Cells.Item(12, 18).Value = 3536.67007016099
$ws.Cells.Item(12, 19).Value = 0.1128671251062779
$ws.Cells.Item(12, 20).Value = 0.1128671251062778
$ws.Cells.Item(13, 7).Value = 9.938311666666667
$ws.Cells.Item(13, 8).Value = 29.814935
$ws.Cells.Item(13, 9).Value = 0.6400052196548212
$ws.Cells.Item(13, 10).Value = 0.640005219654821
$ws.Cells.Item(13, 13).Value = 11.70021233333333
$ws.Cells.Item(13, 14).Value = 35.100637
$ws.Cells.Item(13, 15).Value = 0.05218410806758597
$ws.Cells.Item(13, 16).Value = 0.05218410806758598
$ws.Cells.Item(13, 17).Value = 116.2803567348439
$ws.Cells.Item(13, 18).Value = 1046.523210613595
$ws.Cells.Item(13, 19).Value = 0.03339810154628629
$ws.Cells.Item(13, 20).Value = 0.03339810154628628
$ws.Cells.Item(14, 7).Value = 1.141526666666667
$ws.Cells.Item(14, 8).Value = 3.42458
$ws.Cells.Item(14, 9).Value = 0.0735117844504946
$ws.Cells.Item(14, 10).Value = 0.07351178445049458
$ws.Cells.Item(14, 13).Value = 73.202511
$ws.Cells.Item(14, 14).Value = 219.607533
$ws.Cells.Item(14, 15).Value = 0.3264904632507938
$ws.Cells.Item(14, 16).Value = 0.3264904632507938
$ws.Cells.Item(14, 17).Value = 83.56261837346
$ws.Cells.Item(14, 18).Value = 752.0635653611399
$ws.Cells.Item(14, 19).Value = 0.02400089655963448
$ws.Cells.Item(14, 20).Value = 0.02400089655963448
$ws.Cells.Item(15, 7).Value = 1.141526666666667
$ws.Cells.Item(15, 8).Value = 3.42458
$ws.Cells.Item(15, 9).Value = 0.0735117844504946
$ws.Cells.Item(15, 10).Value = 0.07351178445049458
$ws.Cells.Item(15, 15).Value = 0.4449719839907295
$ws.Cells.Item(15, 16).Value = 0.4449719839907295
$ws.Cells.Item(15, 17).Value = 113.8870143858889
$ws.Cells.Item(15, 18).Value = 1024.983129473
$ws.Cells.Item(15, 19).Value = 0.03271068457363544
$ws.Cells.Item(15, 20).Value = 0.03271068457363543
$ws.Cells.Item(16, 7).Value = 1.141526666666667
$ws.Cells.Item(16, 8).Value = 3.42458
$ws.Cells.Item(16, 9).Value = 0.0735117844504946
$ws.Cells.Item(16, 10).Value = 0.07351178445049458
$ws.Cells.Item(16, 13).Value = 39.54025133333334
$ws.Cells.Item(16, 14).Value = 118.620754
$ws.Cells.Item(16, 15).Value = 0.1763534446908907
$ws.Cells.Item(16, 16).Value = 0.1763534446908907
$ws.Cells.Item(16, 17).Value = 45.13625130370222
$ws.Cells.Item(16, 18).Value = 406.22626173332
$ws.Cells.Item(16, 19).Value = 0.01296405641321898
$ws.Cells.Item(16, 20).Value = 0.01296405641321898
$ws.Cells.Item(17, 7).Value = 1.141526666666667
$ws.Cells.Item(17, 8).Value = 3.42458
$ws.Cells.Item(17, 9).Value = 0.0735117844504946
$ws.Cells.Item(17, 10).Value = 0.07351178445049458
$ws.Cells.Item(17, 13).Value = 11.70021233333333
$ws.Cells.Item(17, 14).Value = 35.100637
$ws.Cells.Item(17, 15).Value = 0.05218410806758597
$ws.Cells.Item(17, 16).Value = 0.05218410806758598
$ws.Cells.Item(17, 17).Value = 13.35610438416222
$ws.Cells.Item(17, 18).Value = 120.20493945746
$ws.Cells.Item(17, 19).Value = 0.003836146904005696
$ws.Cells.Item(17, 20).Value = 0.003836146904005696
